$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row for "e012 Hatches" right after the e011 Deployment row ---
$ws.Rows.Item(13).Insert()

# --- New row 13, column A: the "e012" id tag ---
$ws.Range("A13").Value = "e012"

# --- Update the e011 Deployment body text: add the r4.41 button line ---
$deployText = @'
<Bold>e011 Deployment</Bold> 
<InlineUIContainer><Button Content='r4.4' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>  
<InlineUIContainer><Button Content='r4.41' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>  
<LineBreak/><LineBreak/>
Determine your tank's deployment from the 
<InlineUIContainer><Button Content='Deployment' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>  Table:  
<InlineUIContainer><Image Name='DieRoll' Height='21' Width='21' > </Image></InlineUIContainer>
<LineBreak/><LineBreak/>
'@
$ws.Range("B12").Value = $deployText

# --- New row 13, column B: e012 Hatches body text ---
$hatchesText = @'
<Bold>e012 Hatches</Bold> 
<InlineUIContainer><Button Content='r4.42' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>  
<LineBreak/><LineBreak/>
Left click on hatches on the Tank Card to toggle adding counter. Click image below to continue.
<LineBreak/><LineBreak/>
                                     <InlineUIContainer><Image Name='c15OpenHatch'  Height='100' Width='100'></Image></InlineUIContainer>
'@
$ws.Range("B13").Value = $hatchesText

# --- Row heights (minor re-wrap adjustments from the source edit) ---
$ws.Rows.Item(4).RowHeight = 105
$ws.Rows.Item(5).RowHeight = 99.95
$ws.Rows.Item(6).RowHeight = 114.2
$ws.Rows.Item(8).RowHeight = 99.95
$ws.Rows.Item(9).RowHeight = 99.95
$ws.Rows.Item(10).RowHeight = 156.94999999999999
$ws.Rows.Item(11).RowHeight = 114.2
$ws.Rows.Item(12).RowHeight = 120
$ws.Rows.Item(13).RowHeight = 90
$ws.Rows.Item(14).RowHeight = 28.5
$ws.Rows.Item(15).RowHeight = 28.5
$ws.Rows.Item(16).RowHeight = 128.44999999999999
$ws.Rows.Item(17).RowHeight = 60
$ws.Rows.Item(18).RowHeight = 90

# --- Selection moved to the row that now holds the former B11 content ---
[void]$ws.Range("B15").Select()
